$wb = $excel.ActiveWorkbook

# --- Sheet "AIC": drop the now-unused column C, re-express the filter-step
#     counts in B2:B5 as comma-grouped text, and fix the thousands grouping
#     in the bold AIC figure.
$aic = $wb.Worksheets.Item("AIC")
$aic.Columns.Item(3).Delete()

$aic.Range("B2").NumberFormat = "@"
$aic.Range("B2").Value = "23,212"
$aic.Range("B3").NumberFormat = "@"
$aic.Range("B3").Value = "23,214"
$aic.Range("B4").NumberFormat = "@"
$aic.Range("B4").Value = "22,901"
$aic.Range("B5").NumberFormat = "@"
$aic.Range("B5").Value = "22,642"
$aic.Range("B6").Value = "\textbf{22,341}"

# --- Sheet "Filter": widen the description column, tidy the wording of the
#     depth-bin filter step, and re-express the drift counts in C3:C5 as
#     comma-grouped text.
$flt = $wb.Worksheets.Item("Filter")
$flt.Columns.Item(1).ColumnWidth = 57.2

$flt.Range("A5").Value = "Remove depth bins with little or no sampling (keep 5-39 m)"

$flt.Range("C3").NumberFormat = "@"
$flt.Range("C3").Value = "3,857"
$flt.Range("C4").NumberFormat = "@"
$flt.Range("C4").Value = "3,857"
$flt.Range("C5").NumberFormat = "@"
$flt.Range("C5").Value = "3,848"

# --- Selection / active-sheet bookkeeping: "Filter" becomes the active tab
#     (selected at C6); "Index" keeps its own cursor at F9 but is no longer
#     the active tab; "AIC" keeps its original selection (B7).
$idx = $wb.Worksheets.Item("Index")
$idx.Activate()
$idx.Range("F9").Select()

$flt.Activate()
$flt.Range("C6").Select()
